# Atualização dos dados: 28.12.2025 20:04
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quantidade_pontos")

# New data row (row 15) - mirrors the structure of the existing rows above it.
# Leave most cells' formatting alone so they simply inherit the column's
# default style (same as every other data row in this table).
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 20
$ws.Range("C15").Formula = "=3*23301"
$ws.Range("D15").Formula = "=C15*F15"
$ws.Range("E15").Value = 313
$ws.Range("F15").Value = 1.5
$ws.Range("G15").Value = 1410

# H (Tempo) and L (Data) columns use time/date number formats. Copy the
# format straight from the cell above so the new cells reuse the existing
# style (instead of letting NumberFormat provision a fresh one).
$ws.Range("H14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = 0.037280092592592594

$ws.Range("I15").Value = 7072

# J15 already carried the column's integer-ish style from the template row,
# so just writing the value keeps it consistent with the rest of the table.
$ws.Range("J15").Value = "Vampiro"
$ws.Range("K15").Value = "Normal"

$ws.Range("L14").Copy()
$ws.Range("L15").PasteSpecial(-4122)
$ws.Range("L15").Value = 46020

# Move the active selection to L16, matching where the cursor ended up after the edit.
$ws.Range("L16").Select()
